$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: apply the "A" column style (Mangal 10) used by col A of data rows
function Set-ColAStyle($range) {
    $range.Font.Name = "Mangal"
    $range.Font.Size = 10
}

# Helper: apply the default data style (Arial 10) used by cols B/E/F/G
function Set-DefaultStyle($range) {
    $range.Font.Name = "Arial"
    $range.Font.Size = 10
}

# ---- Row 20 ----
$ws.Range("A20").Value = "com.hamxa.shaynachim"
Set-ColAStyle $ws.Range("A20")

$ws.Range("B20").Value = "bitcoin"
Set-DefaultStyle $ws.Range("B20")

$ws.Range("C20").Value = "echale484@gmail.com"
$ws.Range("D20").Value = "innaplutov1@gmail.com"

$ws.Range("E20").Value = "27/5/2019 15:59"
Set-DefaultStyle $ws.Range("E20")

$ws.Range("F20").Value = "never stop learning. Good guide"
Set-DefaultStyle $ws.Range("F20")

$ws.Range("G20").Value = "yes"
Set-DefaultStyle $ws.Range("G20")

# ---- Row 21 ----
$ws.Range("A21").Value = "com.hamxa.shaynachim"
Set-ColAStyle $ws.Range("A21")

$ws.Range("B21").Value = "bitcoin"
Set-DefaultStyle $ws.Range("B21")

$ws.Range("C21").Value = "leniyadoniv@gmail.com"
$ws.Range("D21").Value = "echale484@gmail.com"

$ws.Range("E21").Value = "27/5/2019 15:59"
Set-DefaultStyle $ws.Range("E21")

$ws.Range("F21").Value = "lets rewrite the rules of bitcoin"
Set-DefaultStyle $ws.Range("F21")

$ws.Range("G21").Value = "yes"
Set-DefaultStyle $ws.Range("G21")

# ---- Row 22 ----
$ws.Range("A22").Value = "com.hamxa.shaynachim"
Set-ColAStyle $ws.Range("A22")

$ws.Range("B22").Value = "bitcoin"
Set-DefaultStyle $ws.Range("B22")

$ws.Range("C22").Value = "rotemzinger3@gmail.com"
$ws.Range("D22").Value = "innaplutov1@gmail.com"

$ws.Range("E22").Value = "27/5/2019 15:59"
Set-DefaultStyle $ws.Range("E22")

$ws.Range("F22").Value = "show me the money please"
Set-DefaultStyle $ws.Range("F22")

$ws.Range("G22").Value = "yes"
Set-DefaultStyle $ws.Range("G22")

# ---- New hyperlinks on the new emails in column C ----
# (Hyperlinks.Add re-fonts the cell with the workbook's built-in "Hyperlink"
#  style, so restore the plain Calibri 11 black look used by every other
#  email cell in columns C/D right after adding each link.)
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:leniyadoniv@gmail.com", "", "", "leniyadoniv@gmail.com")
$ws.Range("C21").Font.Name = "Calibri"
$ws.Range("C21").Font.Size = 11
$ws.Range("C21").Font.Underline = $false
$ws.Range("C21").Font.Color = 0x000000

$ws.Hyperlinks.Add($ws.Range("C22"), "mailto:rotemzinger3@gmail.com", "", "", "rotemzinger3@gmail.com")
$ws.Range("C22").Font.Name = "Calibri"
$ws.Range("C22").Font.Size = 11
$ws.Range("C22").Font.Underline = $false
$ws.Range("C22").Font.Color = 0x000000

# ---- Update the active cell / selection ----
$ws.Range("F23").Select()
